# Commit: "Corregida referencia y diapositivas"
#
# The deck went from 10 slides down to 4:
#   1  OpenAPI Generator           (kept)
#   2  Introducción                (removed)
#   3  Contexto                    (kept)
#   4  Casos de uso                (removed)
#   5  Componentes                 (removed)
#   6  Despliegue                  (kept)
#   7  Desarrollo                  (removed)
#   8  Puntos de variabilidad/extensión (removed)
#   9  Atributos de calidad        (removed)
#   10 FIN                         (kept)
#
# Resulting order: OpenAPI Generator, Contexto, Despliegue, FIN
#
# Delete from the highest index down to the lowest so earlier
# (still-wanted) slide indexes never shift out from under us.

$p = $ppt.ActivePresentation

$p.Slides.Item(9).Delete()   # Atributos de calidad
$p.Slides.Item(8).Delete()   # Puntos de variabilidad/extensión
$p.Slides.Item(7).Delete()   # Desarrollo
$p.Slides.Item(5).Delete()   # Componentes
$p.Slides.Item(4).Delete()   # Casos de uso
$p.Slides.Item(2).Delete()   # Introducción
